# "Generate Report for Archive"
#
# The localization-status report was regenerated: every cell whose status
# was "Ready for handoff" is now "In Translation" (Overview!E2:F2, the
# zh-cn sheet's Status cell C2, and the de-de sheet's Status cell C2).
# Because the new status text is shorter than the old one, the Status
# column(s) that were sized to fit the text shrink accordingly.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns E (zh-cn) and F (de-de).
if ($overview.Range("E2").Text -eq $oldStatus) {
    $overview.Range("E2").Value = $newStatus
}
if ($overview.Range("F2").Text -eq $oldStatus) {
    $overview.Range("F2").Value = $newStatus
}

# Per-language sheets: Status column is column C.
if ($zhcn.Range("C2").Text -eq $oldStatus) {
    $zhcn.Range("C2").Value = $newStatus
}
if ($dede.Range("C2").Text -eq $oldStatus) {
    $dede.Range("C2").Value = $newStatus
}

# The status columns were width-fitted to their contents; re-fit them now
# that the text is shorter so the sheet keeps matching its data.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
